# Updates the crypto price/volume snapshot (GitHub Actions style refresh).
# Values are written with a leading "'" to force text (many Price values,
# e.g. "0.999", look numeric but must stay literal text strings), and the
# range Style is reset to an untouched cell's style afterwards so we don't
# leave a stray number-format/quote-prefix style on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$defaultStyle = $ws.Range("B2").Style
$sub6 = [char]0x2086   # subscript six (₆), used in D51's "0.0₆0221"

$ws.Range("D2").Value = '''63.254.66'
$ws.Range("D2").Style = $defaultStyle
$ws.Range("E2").Value = '''  +6.47%  '
$ws.Range("E2").Style = $defaultStyle
$ws.Range("D3").Value = '''2.436.74'
$ws.Range("D3").Style = $defaultStyle
$ws.Range("E3").Value = '''  +6.15%  '
$ws.Range("E3").Style = $defaultStyle
$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").Style = $defaultStyle
$ws.Range("E4").Value = '''  -0.16%  '
$ws.Range("E4").Style = $defaultStyle
$ws.Range("D5").Value = '''564.87'
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = '''  +5.02%  '
$ws.Range("E5").Style = $defaultStyle
$ws.Range("D6").Value = '''142.35'
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = '''  +11.85%  '
$ws.Range("E6").Style = $defaultStyle
$ws.Range("D7").Value = '''0.999'
$ws.Range("D7").Style = $defaultStyle
$ws.Range("E7").Value = '''  -0.16%  '
$ws.Range("E7").Style = $defaultStyle
$ws.Range("D8").Value = '''0.588'
$ws.Range("D8").Style = $defaultStyle
$ws.Range("E8").Value = '''  +4.26%  '
$ws.Range("E8").Style = $defaultStyle
$ws.Range("D9").Value = '''2.433.97'
$ws.Range("D9").Style = $defaultStyle
$ws.Range("E9").Value = '''  +6.11%  '
$ws.Range("E9").Style = $defaultStyle
$ws.Range("E10").Value = '''  +4.98%  '
$ws.Range("E10").Style = $defaultStyle
$ws.Range("E11").Value = '''  +5.16%  '
$ws.Range("E11").Style = $defaultStyle
$ws.Range("E12").Value = '''  +1.25%  '
$ws.Range("E12").Style = $defaultStyle
$ws.Range("D13").Value = '''0.352'
$ws.Range("D13").Style = $defaultStyle
$ws.Range("E13").Value = '''  +7.20%  '
$ws.Range("E13").Style = $defaultStyle
$ws.Range("D14").Value = '''26.34'
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Value = '''  +15.00%  '
$ws.Range("E14").Style = $defaultStyle
$ws.Range("D15").Value = '''2.867.71'
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = '''  +5.76%  '
$ws.Range("E15").Style = $defaultStyle
$ws.Range("D16").Value = '''63.068.56'
$ws.Range("D16").Style = $defaultStyle
$ws.Range("E16").Value = '''  +6.30%  '
$ws.Range("E16").Style = $defaultStyle
$ws.Range("E17").Value = '''  +9.50%  '
$ws.Range("E17").Style = $defaultStyle
$ws.Range("D18").Value = '''2.433.87'
$ws.Range("D18").Style = $defaultStyle
$ws.Range("E18").Value = '''  +4.54%  '
$ws.Range("E18").Style = $defaultStyle
$ws.Range("D19").Value = '''11.21'
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").Value = '''  +8.55%  '
$ws.Range("E19").Style = $defaultStyle
$ws.Range("D20").Value = '''339.28'
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = '''  +10.30%  '
$ws.Range("E20").Style = $defaultStyle
$ws.Range("D21").Value = '''4.28'
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = '''  +7.26%  '
$ws.Range("E21").Style = $defaultStyle
$ws.Range("E22").Value = '''  +4.70%  '
$ws.Range("E22").Style = $defaultStyle
$ws.Range("E23").Value = '''  +0.09%  '
$ws.Range("E23").Style = $defaultStyle
$ws.Range("D24").Value = '''65.33'
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = '''  +4.07%  '
$ws.Range("E24").Style = $defaultStyle
$ws.Range("E25").Value = '''  +3.59%  '
$ws.Range("E25").Style = $defaultStyle
$ws.Range("E26").Value = '''  +0.00%  '
$ws.Range("E26").Style = $defaultStyle
$ws.Range("E27").Value = '''  +15.30%  '
$ws.Range("E27").Style = $defaultStyle
$ws.Range("E28").Value = '''  +6.67%  '
$ws.Range("E28").Style = $defaultStyle
$ws.Range("D29").Value = '''1.33'
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").Value = '''  +13.45%  '
$ws.Range("E29").Style = $defaultStyle
$ws.Range("D30").Value = '''6.68'
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").Value = '''  +16.74%  '
$ws.Range("E30").Style = $defaultStyle
$ws.Range("E31").Value = '''  +12.32%  '
$ws.Range("E31").Style = $defaultStyle
$ws.Range("E32").Value = '''  +7.51%  '
$ws.Range("E32").Style = $defaultStyle
$ws.Range("D33").Value = '''174.34'
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Value = '''  +1.59%  '
$ws.Range("E33").Style = $defaultStyle
$ws.Range("E34").Value = '''  +12.80%  '
$ws.Range("E34").Style = $defaultStyle
$ws.Range("D35").Value = '''0.399'
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = '''  +6.35%  '
$ws.Range("E35").Style = $defaultStyle
$ws.Range("D36").Value = '''18.71'
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = '''  +6.33%  '
$ws.Range("E36").Style = $defaultStyle
$ws.Range("B37").Value = '''Bittensor'
$ws.Range("B37").Style = $defaultStyle
$ws.Range("C37").Value = '''https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("C37").Style = $defaultStyle
$ws.Range("D37").Value = '''372.59'
$ws.Range("D37").Style = $defaultStyle
$ws.Range("E37").Value = '''  +21.41%  '
$ws.Range("E37").Style = $defaultStyle
$ws.Range("B38").Value = '''NEARProtocol'
$ws.Range("B38").Style = $defaultStyle
$ws.Range("C38").Value = '''https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("C38").Style = $defaultStyle
$ws.Range("D38").Value = '''4.48'
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = '''  +13.71%  '
$ws.Range("E38").Style = $defaultStyle
$ws.Range("E39").Value = '''  +0.01%  '
$ws.Range("E39").Style = $defaultStyle
$ws.Range("E40").Value = '''  -0.16%  '
$ws.Range("E40").Style = $defaultStyle
$ws.Range("E41").Value = '''  +14.84%  '
$ws.Range("E41").Style = $defaultStyle
$ws.Range("D42").Value = '''40.38'
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = '''  +7.44%  '
$ws.Range("E42").Style = $defaultStyle
$ws.Range("D43").Value = '''149.39'
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = '''  +11.02%  '
$ws.Range("E43").Style = $defaultStyle
$ws.Range("D44").Value = '''3.70'
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = '''  +9.68%  '
$ws.Range("E44").Style = $defaultStyle
$ws.Range("D45").Value = '''20.73'
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = '''  +13.53%  '
$ws.Range("E45").Style = $defaultStyle
$ws.Range("D46").Value = '''0.594'
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Value = '''  +5.39%  '
$ws.Range("E46").Style = $defaultStyle
$ws.Range("D47").Value = '''0.0959'
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Value = '''  +3.07%  '
$ws.Range("E47").Style = $defaultStyle
$ws.Range("D48").Value = '''0.0520'
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").Value = '''  +7.07%  '
$ws.Range("E48").Style = $defaultStyle
$ws.Range("E49").Value = '''  +7.55%  '
$ws.Range("E49").Style = $defaultStyle
$ws.Range("D50").Value = '''17.87'
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = '''  +8.22%  '
$ws.Range("E50").Style = $defaultStyle
$ws.Range("D51").Value = '''0.0' + $sub6 + '0221'
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = '''  +0.93%  '
$ws.Range("E51").Style = $defaultStyle
